$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refreshed business listing data (id, business_status, name, rating, user_ratings_total)
# for rows 2-53. Column D (price_level) is untouched since it did not change.
$rows = @(
  @(18, 'OPERATIONAL', 'Arizona Sport Shirts', 5, 38),
  @(47, 'OPERATIONAL', 'Cabela''s', 4.3, 3049),
  @(9, 'OPERATIONAL', 'Champs Sports', 4.2, 38),
  @(19, 'OPERATIONAL', 'Collector''s Den', 3.9, 55),
  @(5, 'OPERATIONAL', 'DICK''S Sporting Goods', 4, 508),
  @(53, 'OPERATIONAL', 'Elevation Golf Indianapolis | Hamilton County Sports', 5, 5),
  @(33, 'OPERATIONAL', 'Finish Line', 4.1, 222),
  @(44, 'OPERATIONAL', 'Genesis Sports Performance', 4.9, 9),
  @(10, 'OPERATIONAL', 'Hibbett Sports', 3.9, 97),
  @(55, 'OPERATIONAL', 'Indiana Sports Corp', 5, 2),
  @(38, 'OPERATIONAL', 'Indy Indoor Sport', 0, 0),
  @(40, 'OPERATIONAL', 'Indy Sport Group', 0, 0),
  @(22, 'OPERATIONAL', 'Indy Sports Performance', 5, 2),
  @(27, 'OPERATIONAL', 'KS&E Sports', 2.8, 34),
  @(57, 'OPERATIONAL', 'Larson and Sons Target Sports', 4.4, 52),
  @(23, 'OPERATIONAL', 'Looking Good Sports Llc', 0, 0),
  @(41, 'OPERATIONAL', 'Marksman Shooting Sports', 4.9, 51),
  @(11, 'OPERATIONAL', 'McCleerey''s Sporting Goods', 4.8, 99),
  @(4, 'OPERATIONAL', 'N & D Sports', 0, 0),
  @(24, 'OPERATIONAL', 'National Institute for Fitness and Sport (NIFS)', 4.3, 52),
  @(49, 'OPERATIONAL', 'Oakley Store', 4.6, 78),
  @(14, 'OPERATIONAL', 'Origyn Sport', 5, 9),
  @(21, 'OPERATIONAL', 'Performance Sports', 3.3, 3),
  @(0, 'OPERATIONAL', 'Play It Again Sports', 4.4, 104),
  @(20, 'OPERATIONAL', 'Ruben Sports', 0, 0),
  @(50, 'OPERATIONAL', 'Sport Clips Haircuts of Carmel', 4, 97),
  @(34, 'OPERATIONAL', 'Sport Clips Haircuts of Carmel - 146th Street', 3.9, 57),
  @(59, 'OPERATIONAL', 'Sport Clips Haircuts of Carmel - Carmel Point', 4.6, 158),
  @(39, 'OPERATIONAL', 'Sport Clips Haircuts of Castleton Crossing', 4.4, 126),
  @(31, 'OPERATIONAL', 'Sport Clips Haircuts of Fishers', 4, 108),
  @(37, 'OPERATIONAL', 'Sport Clips Haircuts of Fishers @ Olio', 4.2, 110),
  @(35, 'OPERATIONAL', 'Sport Clips Haircuts of Geist Oaklandon', 4.6, 127),
  @(54, 'OPERATIONAL', 'Sport Clips Haircuts of German Church Shops', 3.9, 129),
  @(36, 'OPERATIONAL', 'Sport Clips Haircuts of Indianapolis - 96th Street', 4.3, 112),
  @(52, 'OPERATIONAL', 'Sport Clips Haircuts of Medford Place', 4.3, 101),
  @(45, 'OPERATIONAL', 'Sport Clips Haircuts of Noblesville at Stoney Creek Commons', 4.6, 261),
  @(29, 'OPERATIONAL', 'Sport Clips Haircuts of North Keystone', 3.9, 100),
  @(32, 'OPERATIONAL', 'Sport Clips Haircuts of Saxony', 4.1, 80),
  @(42, 'OPERATIONAL', 'Sport Clips Haircuts of Shadeland Place', 4.3, 91),
  @(43, 'OPERATIONAL', 'Sport Clips Haircuts of Westfield', 3.7, 65),
  @(28, 'OPERATIONAL', 'Sport Graphics Inc', 4.3, 10),
  @(2, 'OPERATIONAL', 'Sport Passes', 0, 0),
  @(48, 'OPERATIONAL', 'Sport''n Image', 5, 3),
  @(51, 'OPERATIONAL', 'Sport.ly', 0, 0),
  @(17, 'OPERATIONAL', 'Sports Corporation Inc', 0, 0),
  @(58, 'OPERATIONAL', 'Sports Select', 0, 0),
  @(1, 'OPERATIONAL', 'Sports Spot', 3.3, 13),
  @(30, 'OPERATIONAL', 'Sports Travel & Tickets', 0, 0),
  @(56, 'OPERATIONAL', 'St. Vincent Sports Performance', 5, 1),
  @(16, 'OPERATIONAL', 'Sun Valley Sports', 4.5, 112),
  @(25, 'OPERATIONAL', 'The North Face The Fashion Mall at Keystone', 4.2, 124),
  @(12, 'OPERATIONAL', 'Webster''s Sporting Goods', 4.5, 28)
)

for ($i = 0; $i -lt $rows.Length; $i++) {
  $r = 2 + $i
  $row = $rows[$i]
  $ws.Cells.Item($r, 1).Value = $row[0]
  $ws.Cells.Item($r, 2).Value = $row[1]
  $ws.Cells.Item($r, 3).Value = $row[2]
  $ws.Cells.Item($r, 5).Value = $row[3]
  $ws.Cells.Item($r, 6).Value = $row[4]
}

# Remove the trailing row (Webster's Sporting Goods moved up to row 53; old row 54 dropped)
$ws.Rows.Item(54).Delete()
